$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 42, pushing the existing row 42 (and all
# subsequent rows) down by one. This mirrors the diff, where a brand new
# weekly record is inserted and every later record shifts down one row,
# with the former last record (row 84) becoming row 85.
$ws.Rows.Item(42).EntireRow.Insert()

# Populate the newly inserted row 42 with the new weekly record.
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C42").Value = "Arica y Parinacota"
$ws.Range("D42").Value = 44566
$ws.Range("E42").Value = 15
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100102
$ws.Range("H42").Value = "Cítricos"
$ws.Range("I42").Value = 100102004
$ws.Range("J42").Value = "Mandarina"
$ws.Range("K42").Value = "Murcott"
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 270
$ws.Range("N42").Value = 15000
$ws.Range("O42").Value = 16000
$ws.Range("P42").Value = 15500
$ws.Range("Q42").Value = '$/caja 20 kilos'
$ws.Range("R42").Value = "Región de Coquimbo"
$ws.Range("S42").Value = 775
$ws.Range("T42").Value = 20
